$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update battery capacity value for ID_Battery = 1 from 5000 to 7000
$ws.Range("B2").Value = 7000

# Update selection to B3 to match the saved view state
$ws.Range("B3").Select()
